# Correct False Switch/Case Statement.
#
# 1) Remove the empty <p:timing> stub from the slides that still carry it.
# 2) Rewrite the Switch/Case explanation paragraph on the "Switch/Case
#    Statements" slide, splitting it into new runs and fixing the wording.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Strip the vestigial <p:timing><p:tnLst><p:par>.../p:timing> block.
#    PowerPoint leaves that behind after a slide has been visited once in
#    Slide Show view, even though there are no actual animation effects
#    (MainSequence.Count is 0). Touching the animation sequence (add then
#    immediately delete a no-op effect) makes the host regenerate the
#    timing tree, which collapses back down to nothing and removes the
#    element, exactly like re-saving a "clean" slide in PowerPoint.
# ---------------------------------------------------------------------
$timingSlideIndexes = @(1, 2, 3, 4, 5, 6, 7, 17, 18)
foreach ($idx in $timingSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    $mainSeq = $slide.TimeLine.MainSequence
    if ($mainSeq.Count -eq 0) {
        $shape = $slide.Shapes.Item(1)
        $effect = $mainSeq.AddEffect($shape, 1)
        $effect.Delete()
    }
}

# ---------------------------------------------------------------------
# 2) Fix the Switch/Case body text (slide 20).
#    The original paragraph already carries the "Switch"/"Case"/"Break"/
#    "Default " bold runs; we only need to surgically retype two short
#    phrases in place (same as selecting the words in the UI and typing
#    the replacement), which naturally slices new runs in exactly the
#    two spots the author edited, leaving the rest of the paragraph
#    (including the existing bold runs) untouched.
# ---------------------------------------------------------------------
$rsquo = [char]0x2019

$slide20 = $p.Slides.Item(20)
$body = $slide20.Shapes.Item(2).TextFrame.TextRange

# "...for one variable. The variable being tested..."
#                  ^^^^^^^^^^
#  -> "...for one constant such as an integer. The variable being tested..."
$para1Text = $body.Paragraphs(1, 1).Text
$pos1 = $para1Text.IndexOf("variable. ") + 1
$body.Characters($pos1, 10).Text = "constant such as an integer. "

# "...default case can only appear at the end..."
#                  ^^^^^^^^^
#  -> "...default case should(But doesn't have  to) appear at the end..."
$para1Text = $body.Paragraphs(1, 1).Text
$pos2 = $para1Text.IndexOf("can only ") + 1
$body.Characters($pos2, 9).Text = "should(But doesn${rsquo}t have  to) "
